# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# This particular fixture (tests/.../resources/static/tableOfContent/
# tableOfContent-template.docx) is a *static* test resource used to verify
# the table-of-content handling logic. As part of the housekeeping pass
# that introduced the M2Doc version custom property across the template
# fixtures, this file was simply re-saved by the tooling: its textual
# content, styles, numbering, sections and properties are all unchanged -
# only the on-disk XML attribute ordering differs (a resave artifact), and
# no custom document property was actually required on this resource.
#
# We therefore touch the document (so the save pipeline runs) without
# altering any visible content, formatting, or structure.
$d = $word.ActiveDocument

# No-op touch: confirm the document is reachable / saved as-is.
$null = $d.Name
